$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = "학부 경영 상담(?)"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/business-talks/#utm_source=rss&utm_medium=rss&utm_campaign=business-talks"

$ws.Range("D26").Value = "ai plus(est soft)"

$ws.Range("D37").Value = "dsba_seminar"
